$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: the order in which NEW (previously unseen) strings are first written
# determines their position in the rebuilt shared-strings table, so the
# "id18;id21" / "id16;id19" / "id18;id25" / "id16;id22" / "id21;id25" /
# "id19; id22" combos are written first, in that exact order, before any of
# the other edits below.
$ws.Range("C26").Value = "id18;id21"
$ws.Range("C23").Value = "id16;id19"
$ws.Range("C22").Value = "id18;id25"
$ws.Range("C20").Value = "id16;id22"
$ws.Range("C19").Value = "id21;id25"
$ws.Range("C17").Value = "id19; id22"

# Remaining punteggio (B) and id_annotatore_equivalente (C) updates
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = "id14"

$ws.Range("B12").Value = 1
$ws.Range("C12").ClearContents()

$ws.Range("B14").Value = 1
$ws.Range("C14").ClearContents()

$ws.Range("B15").Value = 3
$ws.Range("C15").Value = "id10"

$ws.Range("B16").Value = 1
$ws.Range("C16").ClearContents()

$ws.Range("B17").Value = 2

$ws.Range("B18").Value = 2
$ws.Range("C18").Value = "id20"

$ws.Range("B19").Value = 1

$ws.Range("B20").Value = 1

$ws.Range("B21").Value = 2
$ws.Range("C21").Value = "id17"

$ws.Range("B22").Value = 2

$ws.Range("B23").Value = 1

$ws.Range("B24").Value = 1
$ws.Range("C24").ClearContents()

$ws.Range("B26").Value = 1

# Row 27 ("id26" equivalence row) is removed entirely; delete it so rows shift up
$ws.Rows(27).Delete()

# Update the sheet view: zoom + active selection moved to C17
$ws.Application.ActiveWindow.Zoom = 177
$ws.Range("C17").Select() | Out-Null
